$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 179 (shifts existing rows 179..242 down to 180..243)
$ws.Rows("179").Insert()

# Populate the newly inserted row 179 with the new data record
$ws.Range("A179").Value2 = 10
$ws.Range("B179").Value2 = "Vega Modelo de Temuco"
$ws.Range("C179").Value2 = "La Araucanía"
$ws.Range("D179").Value2 = 44588
$ws.Range("E179").Value2 = 9
$ws.Range("F179").Value2 = 100112017
$ws.Range("G179").Value2 = "Apio"
$ws.Range("H179").Value2 = "Americana (o)"
$ws.Range("I179").Value2 = "Primera"
$ws.Range("J179").Value2 = 55
$ws.Range("K179").Value2 = 10000
$ws.Range("L179").Value2 = 10000
$ws.Range("M179").Value2 = 10000
$ws.Range("N179").Value2 = "$/docena de matas"
$ws.Range("O179").Value2 = "Provincia del Elquí"
$ws.Range("P179").Value2 = 1667
$ws.Range("Q179").Value2 = 6
$ws.Range("R179").Value2 = "Hortaliza"

# Apply the same date format style used by the other D-column cells
$ws.Range("D179").NumberFormat = $ws.Range("D178").NumberFormat
